$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New loading_percent values for rows 2-25 (data rows 0-23), columns B,D,E,F,G,I,J,M,O.
# Columns C,H,K,L,N remain 0 and column A (index) is untouched.
$newValues = @{
    2 = @{ "B"=10.03591604707802; "D"=8.796942828015835; "E"=12.60781341984324; "F"=30.92064485406651; "G"=3.629112171288599; "I"=23.77891671908844; "J"=9.428389757846121; "M"=27.49615983874695; "O"=23.61489877096263 }
    3 = @{ "B"=9.582141057229396; "D"=8.839853716210484; "E"=12.70081725631764; "F"=30.93139426254933; "G"=3.632184171894186; "I"=23.09541390854298; "J"=9.486217189340508; "M"=26.35820750773802; "O"=23.59387516168431 }
    4 = @{ "B"=9.291282861952386; "D"=8.867523638198167; "E"=12.76096016613953; "F"=30.94998089758091; "G"=3.634168488309737; "I"=22.67141696445289; "J"=9.523568024613585; "M"=25.6308364727196; "O"=23.5893944638361 }
    5 = @{ "B"=9.169793992129312; "D"=8.879132920052644; "E"=12.78623432296682; "F"=30.9605532645494; "G"=3.635001868143477; "I"=22.49785144950431; "J"=9.539253693131721; "M"=25.32753710329781; "O"=23.589682316477 }
    6 = @{ "B"=9.149445712313092; "D"=8.881080808820037; "E"=12.79047735126086; "F"=30.96248936148616; "G"=3.635141747938167; "I"=22.4689929466703; "J"=9.541886397005541; "M"=25.27676850737273; "O"=23.58985757710567 }
    7 = @{ "B"=9.289656246586416; "D"=8.867678853108723; "E"=12.7612979211558; "F"=30.95011136398522; "G"=3.634179627214996; "I"=22.66907895617136; "J"=9.523777683277109; "M"=25.62677351545816; "O"=23.58938979635275 }
    8 = @{ "B"=9.882058472240972; "D"=8.811464642105753; "E"=12.63925144091723; "F"=30.92185469920928; "G"=3.630151090351128; "I"=23.54431890635382; "J"=9.447946461098528; "M"=27.10992048659737; "O"=23.60589799144708 }
    9 = @{ "B"=10.94242800762164; "D"=8.711674771501963; "E"=12.42395918525014; "F"=30.96213332788488; "G"=3.623025473538591; "I"=25.21396690350173; "J"=9.313828580226945; "M"=29.77907426552924; "O"=23.70528414551036 }
    10 = @{ "B"=11.65509735507625; "D"=8.644660532042675; "E"=12.28034813808014; "F"=31.05066493105122; "G"=3.618256711516018; "I"=26.39703232191884; "J"=9.224116515471851; "M"=31.58096553180523; "O"=23.81920081212314 }
    11 = @{ "B"=11.96417293426857; "D"=8.615527975522618; "E"=12.21815965427576; "F"=31.1038039039939; "G"=3.61618735742901; "I"=26.92298357228939; "J"=9.185206140090045; "M"=32.36398883510297; "O"=23.87986518245103 }
    12 = @{ "B"=12.07899061588007; "D"=8.604689649136143; "E"=12.19506090059057; "F"=31.12577657101123; "G"=3.615418031933425; "I"=27.12018549560353; "J"=9.17074402121953; "M"=32.6550835742902; "O"=23.90410158578574 }
    13 = @{ "B"=12.05436209621135; "D"=8.607015284057322; "E"=12.20001560150792; "F"=31.12096211399866; "G"=3.615583085531846; "I"=27.07780495701195; "J"=9.173846595370863; "M"=32.59263422080907; "O"=23.89882576503738 }
    14 = @{ "B"=11.97366384178845; "D"=8.614632426510324; "E"=12.21625027910113; "F"=31.10557454106777; "G"=3.616123778558621; "I"=26.93924796773928; "J"=9.184010878329257; "M"=32.38804673988228; "O"=23.88183384241691 }
    15 = @{ "B"=11.92394311181249; "D"=8.619323322585156; "E"=12.22625315108204; "F"=31.09639009156262; "G"=3.616456827692786; "I"=26.85411623557621; "J"=9.190272244819779; "M"=32.26202131124632; "O"=23.87159016222128 }
    16 = @{ "B"=11.63458994589951; "D"=8.646591548270836; "E"=12.28447541927165; "F"=31.04745113819527; "G"=3.618393953417403; "I"=26.36239600694129; "J"=9.226697562659336; "M"=31.52904221941297; "O"=23.8154136113036 }
    17 = @{ "B"=11.45317114954624; "D"=8.663665462844966; "E"=12.32099663584373; "F"=31.02072488948735; "G"=3.619607865339491; "I"=26.05745607412909; "J"=9.249529370958383; "M"=31.06988018196834; "O"=23.78321146243265 }
    18 = @{ "B"=11.34740283655053; "D"=8.673613292955352; "E"=12.34229836506797; "F"=31.00656382105703; "G"=3.620315490869014; "I"=25.88092239744744; "J"=9.262840517094425; "M"=30.80233827074823; "O"=23.765522594746 }
    19 = @{ "B"=11.31134905337997; "D"=8.677003361695132; "E"=12.34956158078399; "F"=31.00197709094668; "G"=3.620556700368285; "I"=25.82096188238707; "J"=9.267378193555306; "M"=30.71116652069539; "O"=23.75967669890196 }
    20 = @{ "B"=11.4726309365542; "D"=8.661834740564684; "E"=12.31707829452513; "F"=31.02344458335551; "G"=3.619477668480279; "I"=26.09003702078465; "J"=9.247080377392178; "M"=31.11911634872485; "O"=23.78655327541041 }
    21 = @{ "B"=11.99742753477327; "D"=8.612389842093812; "E"=12.21146953705828; "F"=31.11004404878287; "G"=3.615964576634289; "I"=26.98000036464977; "J"=9.181017995799596; "M"=32.44828717016073; "O"=23.88679054555433 }
    22 = @{ "B"=12.32744135890217; "D"=8.581202403762319; "E"=12.1450745812022; "F"=31.17742331587698; "G"=3.613751846088002; "I"=27.55011167773118; "J"=9.139429775603897; "M"=33.28533785679834; "O"=23.95966453591932 }
    23 = @{ "B"=12.15250734220226; "D"=8.597744860361987; "E"=12.18027080407508; "F"=31.14047603565736; "G"=3.614925229392717; "I"=27.24695033046212; "J"=9.161481215645523; "M"=32.84152567077395; "O"=23.920099674761 }
    24 = @{ "B"=11.46383773633131; "D"=8.662661999458447; "E"=12.31884882587007; "F"=31.02221125878461; "G"=3.619536500161242; "I"=26.07531095790086; "J"=9.248186991549851; "M"=31.09686776599778; "O"=23.78503987272012 }
    25 = @{ "B"=10.66694570987877; "D"=8.737559173240928; "E"=12.47963742151085; "F"=30.94091843707681; "G"=3.62487082703486; "I"=24.76892866249995; "J"=9.34855654507594; "M"=29.08411652274701; "O"=23.67120937072243 }
}

foreach ($r in $newValues.Keys) {
    $rowData = $newValues[$r]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$r").Value2 = $rowData[$col]
    }
}